$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5 for the new Bolivian Liga fixture,
# shifting the existing rows 5-8 down to 6-9.
$ws.Rows.Item(5).Insert()

# Rewrite the full data block (rows 2-9, columns A:AO) with the updated odds
# and fixtures for 2025-12-23.
$data = New-Object 'object[,]' 8,41
$data[0,0] = "Australian A-League Men"
$data[0,1] = "'2025-12-23"
$data[0,2] = "05:15:00"
$data[0,3] = "Melbourne City"
$data[0,4] = "Macarthur FC"
$data[0,5] = 3.3
$data[0,6] = 3.6
$data[0,7] = 9.8
$data[0,8] = 12
$data[0,9] = 1.59
$data[0,10] = 1.66
$data[0,11] = 0
$data[0,12] = 0
$data[0,13] = 0
$data[0,14] = 0
$data[0,15] = 1.65
$data[0,16] = 2.44
$data[0,17] = 1.08
$data[0,18] = 11.5
$data[0,19] = 0
$data[0,20] = 0
$data[0,21] = 1.11
$data[0,22] = 1.42
$data[0,23] = 1000
$data[0,24] = 1000
$data[0,25] = 1000
$data[0,26] = 1000
$data[0,27] = 1000
$data[0,28] = 1.73
$data[0,29] = 12
$data[0,30] = 260
$data[0,31] = 1000
$data[0,32] = 4.7
$data[0,33] = 42
$data[0,34] = 1000
$data[0,35] = 1000
$data[0,36] = 36
$data[0,37] = 970
$data[0,38] = 1000
$data[0,39] = 440
$data[0,40] = 1000
$data[1,0] = "Friendly Matches"
$data[1,1] = "'2025-12-23"
$data[1,2] = "09:30:00"
$data[1,3] = "Tombense MG"
$data[1,4] = "Desportiva"
$data[1,5] = 1.41
$data[1,6] = 1.97
$data[1,7] = 2.44
$data[1,8] = 1000
$data[1,9] = 3.2
$data[1,10] = 7
$data[1,11] = 1.31
$data[1,12] = 1.04
$data[1,13] = 1.1
$data[1,14] = 1.25
$data[1,15] = 1.46
$data[1,16] = 1.25
$data[1,17] = 1.18
$data[1,18] = 1.05
$data[1,19] = 1.04
$data[1,20] = 1.04
$data[1,21] = 1.13
$data[1,22] = 2.02
$data[1,23] = 1000
$data[1,24] = 1000
$data[1,25] = 1000
$data[1,26] = 1000
$data[1,27] = 1000
$data[1,28] = 980
$data[1,29] = 1000
$data[1,30] = 1000
$data[1,31] = 1000
$data[1,32] = 1000
$data[1,33] = 60
$data[1,34] = 1000
$data[1,35] = 1000
$data[1,36] = 1000
$data[1,37] = 1000
$data[1,38] = 1000
$data[1,39] = 1000
$data[1,40] = 1000
$data[2,0] = "Algerian Ligue 1"
$data[2,1] = "'2025-12-23"
$data[2,2] = "15:30:00"
$data[2,3] = "MC Alger"
$data[2,4] = "ES Ben Aknoun"
$data[2,5] = 1.35
$data[2,6] = 1.4
$data[2,7] = 14.5
$data[2,8] = 17.5
$data[2,9] = 4.5
$data[2,10] = 5.1
$data[2,11] = 1.5
$data[2,12] = 1.09
$data[2,13] = 3
$data[2,14] = 1.43
$data[2,15] = 1.68
$data[2,16] = 2.28
$data[2,17] = 1.24
$data[2,18] = 4.5
$data[2,19] = 2.82
$data[2,20] = 1.51
$data[2,21] = 1.06
$data[2,22] = 3.5
$data[2,23] = 12
$data[2,24] = 32
$data[2,25] = 1000
$data[2,26] = 1000
$data[2,27] = 5.7
$data[2,28] = 12.5
$data[2,29] = 190
$data[2,30] = 1000
$data[2,31] = 6.4
$data[2,32] = 12
$data[2,33] = 60
$data[2,34] = 1000
$data[2,35] = 11
$data[2,36] = 21
$data[2,37] = 380
$data[2,38] = 1000
$data[2,39] = 10.5
$data[2,40] = 1000
$data[3,0] = "Bolivian Liga de Futbol Profesional"
$data[3,1] = "'2025-12-23"
$data[3,2] = "16:00:00"
$data[3,3] = "Academia de Balompie Boliviano"
$data[3,4] = "San Juan FC"
$data[3,5] = 1.29
$data[3,6] = 1.4
$data[3,7] = 1.09
$data[3,8] = 7.8
$data[3,9] = 1.1
$data[3,10] = 1000
$data[3,11] = 1.21
$data[3,12] = 1.03
$data[3,13] = 1.02
$data[3,14] = 1.13
$data[3,15] = 1.85
$data[3,16] = 1.41
$data[3,17] = 1.7
$data[3,18] = 2.04
$data[3,19] = 1.01
$data[3,20] = 1.01
$data[3,21] = 1.14
$data[3,22] = 3.35
$data[3,23] = 1000
$data[3,24] = 1000
$data[3,25] = 1000
$data[3,26] = 1000
$data[3,27] = 1000
$data[3,28] = 1000
$data[3,29] = 1000
$data[3,30] = 1000
$data[3,31] = 1000
$data[3,32] = 1000
$data[3,33] = 1000
$data[3,34] = 1000
$data[3,35] = 1000
$data[3,36] = 1000
$data[3,37] = 1000
$data[3,38] = 1000
$data[3,39] = 1000
$data[3,40] = 1000
$data[4,0] = "Friendly Matches"
$data[4,1] = "'2025-12-23"
$data[4,2] = "16:00:00"
$data[4,3] = "Serra Branca EC"
$data[4,4] = "Maguary"
$data[4,5] = 2.12
$data[4,6] = 2.5
$data[4,7] = 3.1
$data[4,8] = 3.9
$data[4,9] = 3.25
$data[4,10] = 4.2
$data[4,11] = 1.34
$data[4,12] = 1.06
$data[4,13] = 3.35
$data[4,14] = 1.3
$data[4,15] = 1.9
$data[4,16] = 1.8
$data[4,17] = 1.36
$data[4,18] = 3.05
$data[4,19] = 1.66
$data[4,20] = 2.1
$data[4,21] = 1.36
$data[4,22] = 1.67
$data[4,23] = 16.5
$data[4,24] = 14.5
$data[4,25] = 27
$data[4,26] = 70
$data[4,27] = 11.5
$data[4,28] = 9
$data[4,29] = 16
$data[4,30] = 44
$data[4,31] = 16.5
$data[4,32] = 12.5
$data[4,33] = 18.5
$data[4,34] = 55
$data[4,35] = 34
$data[4,36] = 26
$data[4,37] = 40
$data[4,38] = 200
$data[4,39] = 19
$data[4,40] = 38
$data[5,0] = "Portuguese Primeira Liga"
$data[5,1] = "'2025-12-23"
$data[5,2] = "17:45:00"
$data[5,3] = "Guimaraes"
$data[5,4] = "Sporting Lisbon"
$data[5,5] = 7.8
$data[5,6] = 8
$data[5,7] = 1.51
$data[5,8] = 1.52
$data[5,9] = 4.7
$data[5,10] = 4.8
$data[5,11] = 1.39
$data[5,12] = 1.06
$data[5,13] = 3.95
$data[5,14] = 1.32
$data[5,15] = 2.02
$data[5,16] = 1.95
$data[5,17] = 1.38
$data[5,18] = 3.45
$data[5,19] = 2.08
$data[5,20] = 1.89
$data[5,21] = 2.92
$data[5,22] = 1.14
$data[5,23] = 17
$data[5,24] = 7.6
$data[5,25] = 8.4
$data[5,26] = 13
$data[5,27] = 24
$data[5,28] = 10.5
$data[5,29] = 9.8
$data[5,30] = 16
$data[5,31] = 65
$data[5,32] = 30
$data[5,33] = 28
$data[5,34] = 38
$data[5,35] = 280
$data[5,36] = 130
$data[5,37] = 120
$data[5,38] = 160
$data[5,39] = 180
$data[5,40] = 8.4
$data[6,0] = "Friendly Matches"
$data[6,1] = "'2025-12-23"
$data[6,2] = "18:00:00"
$data[6,3] = "Necaxa"
$data[6,4] = "Atletico San Luis"
$data[6,5] = 1.97
$data[6,6] = 2.34
$data[6,7] = 3.35
$data[6,8] = 4.5
$data[6,9] = 3.15
$data[6,10] = 4
$data[6,11] = 1.36
$data[6,12] = 1.07
$data[6,13] = 3
$data[6,14] = 1.33
$data[6,15] = 1.72
$data[6,16] = 1.85
$data[6,17] = 1.31
$data[6,18] = 3.2
$data[6,19] = 1.76
$data[6,20] = 1.95
$data[6,21] = 1.29
$data[6,22] = 1.74
$data[6,23] = 28
$data[6,24] = 1000
$data[6,25] = 1000
$data[6,26] = 1000
$data[6,27] = 1000
$data[6,28] = 19
$data[6,29] = 1000
$data[6,30] = 1000
$data[6,31] = 1000
$data[6,32] = 1000
$data[6,33] = 60
$data[6,34] = 1000
$data[6,35] = 1000
$data[6,36] = 1000
$data[6,37] = 1000
$data[6,38] = 1000
$data[6,39] = 85
$data[6,40] = 1000
$data[7,0] = "Honduras Liga Nacional"
$data[7,1] = "'2025-12-23"
$data[7,2] = "22:00:00"
$data[7,3] = "Real Espana"
$data[7,4] = "CD Motagua"
$data[7,5] = 1.93
$data[7,6] = 2.16
$data[7,7] = 3.9
$data[7,8] = 4.9
$data[7,9] = 3.3
$data[7,10] = 4
$data[7,11] = 1.42
$data[7,12] = 1.07
$data[7,13] = 3.4
$data[7,14] = 1.36
$data[7,15] = 1.81
$data[7,16] = 1.99
$data[7,17] = 1.3
$data[7,18] = 3.55
$data[7,19] = 1.81
$data[7,20] = 1.94
$data[7,21] = 1.28
$data[7,22] = 1.86
$data[7,23] = 16
$data[7,24] = 17.5
$data[7,25] = 40
$data[7,26] = 120
$data[7,27] = 10
$data[7,28] = 8.6
$data[7,29] = 22
$data[7,30] = 75
$data[7,31] = 13.5
$data[7,32] = 12.5
$data[7,33] = 24
$data[7,34] = 80
$data[7,35] = 30
$data[7,36] = 27
$data[7,37] = 50
$data[7,38] = 1000
$data[7,39] = 21
$data[7,40] = 85

$ws.Range("A2:AO9").Value = $data
